# Apply the edits described by the diff:
# 1. Rename sheet "Cp,Winter" -> "Cp, Winter" (space added after comma)
# 2. Update the stored selection on that sheet from B2 to T13

$wb = $excel.ActiveWorkbook

$originalActiveSheet = $wb.ActiveSheet

$ws = $wb.Worksheets.Item("Cp,Winter")
$ws.Name = "Cp, Winter"

$ws.Activate()
$ws.Range("T13").Select()

$originalActiveSheet.Activate()
